# Horarios actualizados Linea 141 - 694
# Updates schedule data across sheets LP1912, LP1912-215, 6203-6173
$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 11:44:50'
$ws.Range("A3").Value = 'Total filas: 178'
$ws.Range("A37").Value = '06:44:15'
$ws.Range("C37").Value = '11_ETCHEVERRY'
$ws.Range("D37").Value = 32
$ws.Range("A38").Value = '06:59:44'
$ws.Range("C38").Value = '16_SANTA ANA'
$ws.Range("D38").Value = 17
$ws.Range("A47").Value = '06:44:15'
$ws.Range("C47").Value = '11_ETCHEVERRY'
$ws.Range("D47").Value = 48
$ws.Range("A48").Value = '06:01:37'
$ws.Range("C48").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D48").Value = 91
$ws.Range("A94").Value = '08:20:43'
$ws.Range("C94").Value = '16_SANTA ANA'
$ws.Range("D94").Value = 62
$ws.Range("A95").Value = '07:31:43'
$ws.Range("C95").Value = '17_ROMERO'
$ws.Range("D95").Value = 111
$ws.Range("A151").Value = '11:44:50'
$ws.Range("B151").Value = '11:44'
$ws.Range("C151").Value = '16_SANTA ANA'
$ws.Range("D151").Value = 0
$ws.Range("A152").Value = '11:44:50'
$ws.Range("B152").Value = '11:44'
$ws.Range("C152").Value = '16_P MOR-SANTA ANA'
$ws.Range("D152").Value = 0
$ws.Range("A153").Value = '10:30:21'
$ws.Range("B153").Value = '11:48'
$ws.Range("C153").Value = '10_OLMOS'
$ws.Range("D153").Value = 78
$ws.Range("B154").Value = '11:51'
$ws.Range("C154").Value = '215B_EL PATO'
$ws.Range("D154").Value = 81
$ws.Range("A155").Value = '11:03:46'
$ws.Range("B155").Value = '11:52'
$ws.Range("C155").Value = '15_ABASTO'
$ws.Range("D155").Value = 49
$ws.Range("B156").Value = '11:54'
$ws.Range("C156").Value = '15_ABASTO'
$ws.Range("D156").Value = 84
$ws.Range("A157").Value = '11:44:50'
$ws.Range("B157").Value = '11:58'
$ws.Range("C157").Value = '225_GOMEZ'
$ws.Range("D157").Value = 14
$ws.Range("A158").Value = '10:30:21'
$ws.Range("B158").Value = '11:59'
$ws.Range("C158").Value = '225_GOMEZ'
$ws.Range("D158").Value = 89
$ws.Range("A159").Value = '10:30:21'
$ws.Range("B159").Value = '12:02'
$ws.Range("C159").Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Range("D159").Value = 92
$ws.Range("A160").Value = '11:44:50'
$ws.Range("B160").Value = '12:04'
$ws.Range("C160").Value = '23_HERNANDEZ'
$ws.Range("D160").Value = 20
$ws.Range("B161").Value = '12:06'
$ws.Range("C161").Value = '16_P MOR-SANTA ANA'
$ws.Range("D161").Value = 96
$ws.Range("A162").Value = '11:03:46'
$ws.Range("B162").Value = '12:06'
$ws.Range("C162").Value = '14_ABASTO'
$ws.Range("D162").Value = 63
$ws.Range("B163").Value = '12:09'
$ws.Range("C163").Value = '10_OLMOS'
$ws.Range("D163").Value = 66
$ws.Range("B164").Value = '12:14'
$ws.Range("C164").Value = '17_ROMERO'
$ws.Range("D164").Value = 104
$ws.Range("A165").Value = '11:44:50'
$ws.Range("B165").Value = '12:16'
$ws.Range("C165").Value = '16_SANTA ANA'
$ws.Range("D165").Value = 32
$ws.Range("A166").Value = '10:30:21'
$ws.Range("B166").Value = '12:17'
$ws.Range("C166").Value = '14_ABASTO'
$ws.Range("D166").Value = 107
$ws.Range("B167").Value = '12:20'
$ws.Range("C167").Value = '14_ABASTO'
$ws.Range("D167").Value = 77
$ws.Range("A168").Value = '10:30:21'
$ws.Range("B168").Value = '12:20'
$ws.Range("C168").Value = '215A_EL PATO'
$ws.Range("D168").Value = 110
$ws.Range("A169").Value = '10:30:21'
$ws.Range("B169").Value = '12:21'
$ws.Range("C169").Value = '26_HERNANDEZ'
$ws.Range("D169").Value = 111
$ws.Range("E169").Value = 'LP1912'
$ws.Range("A170").Value = '11:03:46'
$ws.Range("B170").Value = '12:31'
$ws.Range("C170").Value = '17_ROMERO'
$ws.Range("D170").Value = 88
$ws.Range("E170").Value = 'LP1912'
$ws.Range("A171").Value = '11:44:50'
$ws.Range("B171").Value = '12:34'
$ws.Range("C171").Value = '23_HERNANDEZ'
$ws.Range("D171").Value = 50
$ws.Range("E171").Value = 'LP1912'
$ws.Range("A172").Value = '11:03:46'
$ws.Range("B172").Value = '12:36'
$ws.Range("C172").Value = '27_EL RETIRO'
$ws.Range("D172").Value = 93
$ws.Range("E172").Value = 'LP1912'
$ws.Range("A173").Value = '11:44:50'
$ws.Range("B173").Value = '12:37'
$ws.Range("C173").Value = '11_ETCHEVERRY'
$ws.Range("D173").Value = 53
$ws.Range("E173").Value = 'LP1912'
$ws.Range("A174").Value = '11:03:46'
$ws.Range("B174").Value = '12:38'
$ws.Range("C174").Value = '17_179 Y 38'
$ws.Range("D174").Value = 95
$ws.Range("E174").Value = 'LP1912'
$ws.Range("A175").Value = '11:44:50'
$ws.Range("B175").Value = '12:41'
$ws.Range("C175").Value = '10_OLMOS'
$ws.Range("D175").Value = 57
$ws.Range("E175").Value = 'LP1912'
$ws.Range("A176").Value = '11:03:46'
$ws.Range("B176").Value = '12:48'
$ws.Range("C176").Value = '11_ETCHEVERRY'
$ws.Range("D176").Value = 105
$ws.Range("E176").Value = 'LP1912'
$ws.Range("A177").Value = '11:44:50'
$ws.Range("B177").Value = '13:02'
$ws.Range("C177").Value = '15_ABASTO'
$ws.Range("D177").Value = 78
$ws.Range("E177").Value = 'LP1912'
$ws.Range("A178").Value = '11:44:50'
$ws.Range("B178").Value = '13:10'
$ws.Range("C178").Value = '10_OLMOS'
$ws.Range("D178").Value = 86
$ws.Range("E178").Value = 'LP1912'
$ws.Range("A179").Value = '11:44:50'
$ws.Range("B179").Value = '13:13'
$ws.Range("C179").Value = '215D_EL PATO'
$ws.Range("D179").Value = 89
$ws.Range("E179").Value = 'LP1912'
$ws.Range("A180").Value = '11:44:50'
$ws.Range("B180").Value = '13:19'
$ws.Range("C180").Value = '10_OLMOS'
$ws.Range("D180").Value = 95
$ws.Range("E180").Value = 'LP1912'
$ws.Range("A181").Value = '11:44:50'
$ws.Range("B181").Value = '13:21'
$ws.Range("C181").Value = '26_HERNANDEZ'
$ws.Range("D181").Value = 97
$ws.Range("E181").Value = 'LP1912'
$ws.Range("A182").Value = '11:44:50'
$ws.Range("B182").Value = '13:26'
$ws.Range("C182").Value = '14_ABASTO'
$ws.Range("D182").Value = 102
$ws.Range("E182").Value = 'LP1912'
$ws.Range("A183").Value = '11:44:50'
$ws.Range("B183").Value = '13:26'
$ws.Range("C183").Value = '15_ABASTO'
$ws.Range("D183").Value = 102
$ws.Range("E183").Value = 'LP1912'

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 11:44:50'
$ws.Range("A3").Value = 'Total filas: 19'
$ws.Range("A24").Value = '11:44:50'
$ws.Range("B24").Value = '13:13'
$ws.Range("C24").Value = '215D_EL PATO'
$ws.Range("D24").Value = 89
$ws.Range("E24").Value = 'LP1912'

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 11:44:50'
$ws.Range("A3").Value = 'Total filas: 28'
$ws.Range("A33").Value = '11:44:50'
$ws.Range("B33").Value = '13:30'
$ws.Range("C33").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D33").Value = 106
$ws.Range("E33").Value = 'L6173'
